$wb = $excel.ActiveWorkbook

# Sheet 1: دانشجویان (Students) - National ID column (A) changes from
# numeric values to text values for rows 2-4 (values unchanged).
$ws1 = $wb.Worksheets.Item("دانشجویان")
$ws1.Range("A2").Value = "'1234567890"
$ws1.Range("A3").Value = "'9876543210"
$ws1.Range("A4").Value = "'1122334455"

# Sheet 2: پرداخت‌ها (Payments) - add a new column E "کد ملی" (National ID)
# with each payment row's person's national ID looked up from Sheet 1.
$ws2 = $wb.Worksheets.Item("پرداخت‌ها")

# Match the bold/centered/bordered header style used by A1:D1.
$ws2.Range("D1").Copy($ws2.Range("E1"))
$ws2.Range("E1").Value = "کد ملی"

$ws2.Range("E2").Value = "'1122334455"
$ws2.Range("E3").Value = "'1122334455"
$ws2.Range("E4").Value = "'1122334455"
$ws2.Range("E5").Value = "'1122334455"
$ws2.Range("E6").Value = "'1234567890"
$ws2.Range("E7").Value = "'1234567890"
